# "Caminhonete 100% atualizada bomba patch"
#
# What changed, reconstructed from the OOXML diff:
#   1. The row holding "Caminhao" (row 4) was deleted entirely, shifting every
#      row below it up by one (Vendedor, formerly row 16, falls off the used
#      range and row 16 becomes blank).
#   2. "Automovel" (now row 3) was marked 100% (B3: 0 -> 1).
#   3. "Caminhonete" (now row 4, after the shift) was marked 100% (B4: 0 -> 1).
#   4. "Carro" (now row 5, after the shift) keeps its 100% value but is
#      re-styled as bold + underlined + centered percentage, to call it out
#      (new font + new cell style picked up automatically by Excel).
#   5. The active selection ends up on Q24.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Delete the "Caminhao" row (row 4) - Excel shifts everything below up.
$ws.Rows(4).Delete()

# 2. Automovel -> 100%
$ws.Range("B3").Value = 1

# 3. Caminhonete (shifted into row 4) -> 100%
$ws.Range("B4").Value = 1

# 4. Carro (shifted into row 5) -> emphasize with bold/underline/percent,
#    centered, while keeping its value at 100%.
$ws.Range("B5").Value = 1
$ws.Range("B5").Font.Bold = $true
$ws.Range("B5").Font.Underline = $true
$ws.Range("B5").NumberFormat = "0%"
$ws.Range("B5").HorizontalAlignment = -4108

# 5. Leave the selection where the editor ended up.
$ws.Range("Q24").Select()
